$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.670.82'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '''1.803.32'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").Value = '''0.9997'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''306.55'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''0.9992'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '''0.4300'
$ws.Range("E7").Value = '  +2.40%  '
$ws.Range("D8").Value = '''0.3651'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("D9").Value = '''0.07186'
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("D10").Value = '''0.8593'
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("D11").Value = '''20.76'
$ws.Range("E11").Value = '  +3.04%  '
$ws.Range("D12").Value = '''1.916.00'
$ws.Range("E12").Value = '  +6.11%  '
$ws.Range("D13").Value = '''6.584'
$ws.Range("E13").Value = '  +3.44%  '
$ws.Range("D14").Value = '''5.320'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '''0.06866'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("D16").Value = '''1.005'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '''80.14'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '''0.000008796'
$ws.Range("E18").Value = '  +1.57%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '''15.19'
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '''26.675.30'
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("D22").Value = '''5.171'
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("D23").Value = '''11.08'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").Value = '''2.131.15'
$ws.Range("E24").Value = '  +5.48%  '
$ws.Range("D25").Value = '''152.12'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '''1.841'
$ws.Range("E26").Value = '  -4.25%  '
$ws.Range("D27").Value = '''18.24'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '''5.192'
$ws.Range("E28").Value = '  +3.59%  '
$ws.Range("D29").Value = '''1.903'
$ws.Range("E29").Value = '  +15.89%  '
$ws.Range("D30").Value = '''115.17'
$ws.Range("E30").Value = '  +1.97%  '
$ws.Range("D31").Value = '''0.08940'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = '''0.7497'
$ws.Range("E32").Value = '  +3.57%  '
$ws.Range("D33").Value = '''1.158'
$ws.Range("E33").Value = '  +6.60%  '
$ws.Range("D34").Value = '''4.397'
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("D35").Value = '''2.757'
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("D36").Value = '''1.001'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '''1.118'
$ws.Range("E37").Value = '  +3.62%  '
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = '''0.01911'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '''0.5027'
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").Value = '''0.1633'
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = '''2.634'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = '''6.475'
$ws.Range("E43").Value = '  +9.63%  '
$ws.Range("D44").Value = '''8.227'
$ws.Range("E44").Value = '  +2.48%  '
$ws.Range("D45").Value = '''106.14'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("D46").Value = '''10.28'
$ws.Range("E46").Value = '  +1.28%  '
$ws.Range("D47").Value = '''0.9994'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '''1.647'
$ws.Range("E48").Value = '  +2.87%  '
$ws.Range("D49").Value = '''0.4538'
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("D50").Value = '''0.06233'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").Value = '''1.796'
$ws.Range("E51").Value = '  +5.59%  '
